{"js": "// Word JS API (Office.js) edit script.\n// Applies the two visible-text content changes from the commit:\n//   1. \"...120 were of the A5 variant.\" -> \"...120 were of the very modern A5 variant.\"\n//   2. \"...air defence systems the OTO-Breda...\" -> \"...air defence systems, the OTO-Breda...\"\n\nconst body = context.document.body;\n\n// 1. Insert \"very modern \" immediately before \"A5 variant\"\nconst a5Results = body.search(\"A5 variant\", { matchCase: true, matchWholeWord: false });\na5Results.load(\"items\");\nawait context.sync();\n\nif (a5Results.items.length > 0) {\n  a5Results.items[0].insertText(\"very modern \", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// 2. Insert \",\" right after \"air defence systems\" (before \" the OTO-Breda...\")\nconst adResults = body.search(\"air defence systems\", { matchCase: true, matchWholeWord: false });\nadResults.load(\"items\");\nawait context.sync();\n\nif (adResults.items.length > 0) {\n  adResults.items[0].insertText(\",\", Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# Applies the two visible-text content changes from the commit:\n#   1. \"...120 were of the A5 variant.\" -> \"...120 were of the very modern A5 variant.\"\n#   2. \"...air defence systems the OTO-Breda...\" -> \"...air defence systems, the OTO-Breda...\"\n#\n# wdReplaceAll = 2, wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n# 1. Insert \"very modern \" immediately before \"A5 variant\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"A5 variant\",   # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    \"very modern A5 variant\",  # ReplaceWith\n    2               # Replace (wdReplaceAll)\n)\n\n# 2. Insert \",\" right after \"air defence systems\" (before \" the OTO-Breda...\")\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"air defence systems the\",  # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    \"air defence systems, the\",  # ReplaceWith\n    2               # Replace (wdReplaceAll)\n)\n"}
